$d = $word.ActiveDocument

# Locate the paragraph that begins the "Creating a county" section (the
# paragraph that currently comes right after "Map Notes").
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim().StartsWith("Creating a county")) {
        $targetIndex = $i
        break
    }
}

# Insert three new paragraphs immediately before it:
#   1) a paragraph with the wikimedia URL
#   2) an empty paragraph
#   3) an empty paragraph
# Each InsertParagraphBefore() call pushes "Creating a county" (and
# everything after it) one slot further down, so re-resolve the
# paragraph by its (growing) index each time rather than caching a
# stale reference.
$d.Paragraphs.Item($targetIndex).Range.InsertParagraphBefore()
$d.Paragraphs.Item($targetIndex + 1).Range.InsertParagraphBefore()
$d.Paragraphs.Item($targetIndex + 2).Range.InsertParagraphBefore()

$d.Paragraphs.Item($targetIndex).Range.Text = "https://commons.wikimedia.org/wiki/File:USA_Oregon_relief_location_map.svg"
